$d = $word.ActiveDocument

# Locate the last paragraph in the document (currently "ASPX") and insert
# four new list paragraphs after it, matching the ListParagraph / numId=1
# numbered-list style already used throughout the document.

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)

# 1) IDE: Integrated Development Environment (Visual Studio)
$lastPara.Range.InsertParagraphAfter()
$p1 = $d.Paragraphs.Item($d.Paragraphs.Count)
$r1 = $p1.Range
$r1.InsertAfter("IDE: Integrated Development Environment")
$r1.Collapse(0)
$r1.InsertAfter(" (Visual Studio)")

# 2) Framework: is simply a collection of packages and utility functions,
#    and functionalities (.NET framework)
$p1.Range.InsertParagraphAfter()
$p2 = $d.Paragraphs.Item($d.Paragraphs.Count)
$r2 = $p2.Range
$r2.InsertAfter("Framework:")
$r2.Collapse(0)
$r2.InsertAfter(" is simply a collection of packages and utility functions, and functionalities")
$r2.Collapse(0)
$r2.InsertAfter(" (.NET framework)")

# 3) CLI: (Command line input)
$p2.Range.InsertParagraphAfter()
$p3 = $d.Paragraphs.Item($d.Paragraphs.Count)
$r3 = $p3.Range
$r3.InsertAfter("CLI: (Command line input)")

# 4) SDK: Software development kit
$p3.Range.InsertParagraphAfter()
$p4 = $d.Paragraphs.Item($d.Paragraphs.Count)
$r4 = $p4.Range
$r4.InsertAfter("SDK: Software development kit")
